$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report header strings (new week) ---
$ws.Range("A8").Value = "Volume 31   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/12/2024  Through  2/18/2024"

# --- Simple numeric updates (style/type unchanged) ---
$ws.Range("D15").Value2 = 3
$ws.Range("E15").Value2 = -66.666666666666
$ws.Range("F15").Value2 = 2
$ws.Range("G15").Value2 = 5
$ws.Range("H15").Value2 = -60
$ws.Range("I15").Value2 = 3
$ws.Range("J15").Value2 = 6
$ws.Range("K15").Value2 = -50
$ws.Range("L15").Value2 = -66.666666666666
$ws.Range("N15").Value2 = -76.923076923076
$ws.Range("D16").Value2 = 5
$ws.Range("E16").Value2 = -100
$ws.Range("F16").Value2 = 9
$ws.Range("G16").Value2 = 13
$ws.Range("H16").Value2 = -30.76923076923
$ws.Range("I16").Value2 = 17
$ws.Range("J16").Value2 = 17
$ws.Range("K16").Value2 = 0
$ws.Range("L16").Value2 = 13.333333333333
$ws.Range("M16").Value2 = -58.536585365853
$ws.Range("N16").Value2 = -88.194444444444
$ws.Range("C17").Value2 = 12
$ws.Range("D17").Value2 = 8
$ws.Range("E17").Value2 = 50
$ws.Range("F17").Value2 = 29
$ws.Range("G17").Value2 = 32
$ws.Range("H17").Value2 = -9.375
$ws.Range("I17").Value2 = 54
$ws.Range("J17").Value2 = 52
$ws.Range("K17").Value2 = 3.846153846153
$ws.Range("L17").Value2 = 25.581395348837
$ws.Range("M17").Value2 = 58.823529411764
$ws.Range("N17").Value2 = -41.304347826087
$ws.Range("C18").Value2 = 1
$ws.Range("D18").Value2 = 5
$ws.Range("E18").Value2 = -80
$ws.Range("G18").Value2 = 10
$ws.Range("H18").Value2 = -20
$ws.Range("I18").Value2 = 11
$ws.Range("J18").Value2 = 13
$ws.Range("K18").Value2 = -15.384615384615
$ws.Range("L18").Value2 = -47.619047619047
$ws.Range("M18").Value2 = -77.083333333333
$ws.Range("N18").Value2 = -90.350877192982
$ws.Range("C19").Value2 = 5
$ws.Range("D19").Value2 = 9
$ws.Range("E19").Value2 = -44.444444444444
$ws.Range("F19").Value2 = 22
$ws.Range("G19").Value2 = 31
$ws.Range("H19").Value2 = -29.032258064516
$ws.Range("I19").Value2 = 36
$ws.Range("J19").Value2 = 52
$ws.Range("K19").Value2 = -30.76923076923
$ws.Range("L19").Value2 = -26.530612244898
$ws.Range("M19").Value2 = -59.090909090909
$ws.Range("N19").Value2 = -84.615384615384
$ws.Range("C20").Value2 = 5
$ws.Range("D20").Value2 = 7
$ws.Range("E20").Value2 = -28.571428571428
$ws.Range("F20").Value2 = 18
$ws.Range("G20").Value2 = 16
$ws.Range("H20").Value2 = 12.5
$ws.Range("I20").Value2 = 30
$ws.Range("J20").Value2 = 24
$ws.Range("K20").Value2 = 25
$ws.Range("L20").Value2 = -9.090909090909
$ws.Range("M20").Value2 = -11.764705882352
$ws.Range("N20").Value2 = -86.486486486486
$ws.Range("C21").Value2 = 24
$ws.Range("D21").Value2 = 37
$ws.Range("E21").Value2 = -35.135135135135
$ws.Range("F21").Value2 = 88
$ws.Range("G21").Value2 = 107
$ws.Range("H21").Value2 = -17.757009345794
$ws.Range("I21").Value2 = 151
$ws.Range("J21").Value2 = 164
$ws.Range("K21").Value2 = -7.926829268292
$ws.Range("L21").Value2 = -12.209302325581
$ws.Range("M21").Value2 = -38.866396761133
$ws.Range("N21").Value2 = -81.696969696969
$ws.Range("I23").Value2 = 2
$ws.Range("L23").Value2 = 100
$ws.Range("M23").Value2 = -33.333333333333
$ws.Range("C24").Value2 = 25
$ws.Range("D24").Value2 = 17
$ws.Range("E24").Value2 = 47.058823529411
$ws.Range("F24").Value2 = 90
$ws.Range("G24").Value2 = 88
$ws.Range("H24").Value2 = 2.272727272727
$ws.Range("I24").Value2 = 156
$ws.Range("J24").Value2 = 152
$ws.Range("K24").Value2 = 2.631578947368
$ws.Range("L24").Value2 = -18.75
$ws.Range("M24").Value2 = 24.8
$ws.Range("C25").Value2 = 13
$ws.Range("D25").Value2 = 15
$ws.Range("E25").Value2 = -13.333333333333
$ws.Range("F25").Value2 = 60
$ws.Range("H25").Value2 = 42.857142857142
$ws.Range("I25").Value2 = 82
$ws.Range("J25").Value2 = 81
$ws.Range("K25").Value2 = 1.234567901234
$ws.Range("L25").Value2 = 46.428571428571
$ws.Range("M25").Value2 = -7.865168539325
$ws.Range("C26").Value2 = 2
$ws.Range("D26").Value2 = 5
$ws.Range("E26").Value2 = -60
$ws.Range("F26").Value2 = 4
$ws.Range("G26").Value2 = 7
$ws.Range("H26").Value2 = -42.857142857142
$ws.Range("I26").Value2 = 5
$ws.Range("J26").Value2 = 9
$ws.Range("K26").Value2 = -44.444444444444
$ws.Range("L26").Value2 = -58.333333333333
$ws.Range("F27").Value2 = 2
$ws.Range("G27").Value2 = 4
$ws.Range("H27").Value2 = -50
$ws.Range("L27").Value2 = 0
$ws.Range("G28").Value2 = 1
$ws.Range("H28").Value2 = 100
$ws.Range("M28").Value2 = -66.666666666666
$ws.Range("N28").Value2 = -91.666666666666
$ws.Range("G29").Value2 = 1
$ws.Range("H29").Value2 = 100
$ws.Range("M29").Value2 = -60
$ws.Range("N29").Value2 = -89.473684210526

# --- Cells that change "kind" (text <-> number) or style: copy formatting from a donor cell first ---
$ws.Range("C36").Copy($ws.Range("C15"))
$ws.Range("C15").Value2 = 1.0
$ws.Range("C22").Copy($ws.Range("C16"))
$ws.Range("C36").Copy($ws.Range("D23"))
$ws.Range("D23").Value2 = 1.0
$ws.Range("K36").Copy($ws.Range("E23"))
$ws.Range("E23").Value2 = -100.0
$ws.Range("C36").Copy($ws.Range("F23"))
$ws.Range("F23").Value2 = 1.0
$ws.Range("C36").Copy($ws.Range("G23"))
$ws.Range("G23").Value2 = 1.0
$ws.Range("K36").Copy($ws.Range("H23"))
$ws.Range("H23").Value2 = 0.0
$ws.Range("C36").Copy($ws.Range("J23"))
$ws.Range("J23").Value2 = 1.0
$ws.Range("K36").Copy($ws.Range("K23"))
$ws.Range("K23").Value2 = 100.0
$ws.Range("C22").Copy($ws.Range("C27"))
$ws.Range("C22").Copy($ws.Range("D27"))
$ws.Range("E22").Copy($ws.Range("E27"))
$ws.Range("C22").Copy($ws.Range("D28"))
$ws.Range("E22").Copy($ws.Range("E28"))
$ws.Range("C22").Copy($ws.Range("D29"))
$ws.Range("E22").Copy($ws.Range("E29"))
